$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44895
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("P2").Value = 1385

$ws.Range("D3").Value = 44797
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12500
$ws.Range("P3").Value = 962

$ws.Range("D4").Value = 44915
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("P4").Value = 1385

$ws.Range("D5").Value = 44868
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("P5").Value = 1385

$ws.Range("D6").Value = 44880
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("P6").Value = 1308

$ws.Range("D7").Value = 44839
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 1192

$ws.Range("D8").Value = 44804
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12500
$ws.Range("P8").Value = 962

$ws.Range("D9").Value = 44922
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 17000
$ws.Range("P9").Value = 1308

$ws.Range("D10").Value = 44810
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("P10").Value = 962

$ws.Range("D11").Value = 44874
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 17000
$ws.Range("M11").Value = 17000
$ws.Range("P11").Value = 1308

$ws.Range("D12").Value = 44959
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 19000
$ws.Range("L12").Value = 19000
$ws.Range("M12").Value = 19000
$ws.Range("P12").Value = 1462

$ws.Range("D13").Value = 44943
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 17000
$ws.Range("P13").Value = 1308

$ws.Range("D14").Value = 44930
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 17000
$ws.Range("P14").Value = 1308

$ws.Range("D15").Value = 44859
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 13000
$ws.Range("P15").Value = 1000

$ws.Range("D16").Value = 44841
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("P16").Value = 1385

$ws.Range("D17").Value = 44894
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 18000
$ws.Range("P17").Value = 1385

$ws.Range("D18").Value = 44832
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 1346

$ws.Range("D19").Value = 44846
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 18000
$ws.Range("P19").Value = 1385

